# automation script for admin
#
# Updates the sample row (row 2) on Sheet1 with freshly generated
# test-data values (customer name, site name, dispatcher/driver notes,
# customer ID) and records a Payment status. Also nudges the sheet view
# (scroll position / selection) and widens the new Payment column to
# match where the analyst left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- refreshed sample data on row 2 -----------------------------------
$ws.Range("C2").Value = "ryLnl_0405417"   # CustomerName
$ws.Range("D2").Value = "UnWKKNELUy"      # SiteName
$ws.Range("K2").Value = "FRhWixaQjb"      # Dispatcher Note
$ws.Range("M2").Value = "FRhWixaQjb"      # Driver Note
$ws.Range("N2").Value = "UCN 10540"       # CustomerID
$ws.Range("R2").Value = "Pass"            # Payment

# --- widen the new Payment column (P) ---------------------------------
$ws.Columns.Item(16).ColumnWidth = 13.67

# --- restore the view state the workbook was saved with ---------------
$ws.Range("N7").Select()
